$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @("IMX-USD", "TAO-USD", "GRT-USD", "MNT-USD")

$startRow = 346
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newValues[$i]
}
